# Insert a new row before row 660, pushing the existing row 660 (and all
# rows below it) down by one. Then populate the newly inserted row 660
# with the "new" data record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(660).Insert()

$row = 660

$ws.Cells.Item($row, 1).Value = 6
$ws.Cells.Item($row, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item($row, 3).Value = "Metropolitana"
$ws.Cells.Item($row, 4).Value = 45142
$ws.Cells.Item($row, 5).Value = 13
$ws.Cells.Item($row, 6).Value = 100112039
$ws.Cells.Item($row, 7).Value = "Ciboulette"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 660
$ws.Cells.Item($row, 11).Value = 1200
$ws.Cells.Item($row, 12).Value = 1300
$ws.Cells.Item($row, 13).Value = 1255
$ws.Cells.Item($row, 14).Value = "`$/docena de atados"
$ws.Cells.Item($row, 15).Value = "Región Metropolitana"
$ws.Cells.Item($row, 16).Value = 418
$ws.Cells.Item($row, 17).Value = 3
$ws.Cells.Item($row, 18).Value = "Hortaliza"
